# build: change name to sheets
#
# - Rename sheets "General" -> "general" and "Data" -> "data"
#   (the workbook's _xlnm._FilterDatabase defined name, which is scoped to
#   the first sheet and refers to it by name, follows automatically).
# - Make "general" (sheet 1) the active/selected tab instead of "data"
#   (sheet 2).
# - Update the lingering selection on the "data" sheet from C19 to B19.

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item(1)
$wsData    = $wb.Worksheets.Item(2)

$wsGeneral.Name = "general"
$wsData.Name    = "data"

# Move the saved selection on the (now inactive) "data" sheet to B19.
$wsData.Range("B19").Select()

# Activate "general" so it becomes the selected/visible tab on open.
$wsGeneral.Activate()
